$wb = $excel.ActiveWorkbook

# Sheet "Filenames to Set Mapping": rename the flag value in column C, row 4
# from "sixmix_flag" to "sixmix"
$ws1 = $wb.Worksheets.Item("Filenames to Set Mapping")
$ws1.Range("C4").Value = "sixmix"

# Sheet "Standards per Set": row 8 (set 2) is changed to reference the same
# standard/adduct as row 3 (Sulfamethizole, [M+H]+) instead of
# Sulfadiamethoxine / [M+Na]+, and the mz value changes from 198 to 144.
$ws2 = $wb.Worksheets.Item("Standards per Set")
$ws2.Range("B8").Value = "Sulfamethizole "
$ws2.Range("C8").Value = "C9H10N4O2S2"
$ws2.Range("D8").Value = "[M+H]+"
$ws2.Range("F8").Value = 144
